# fix if the csv file has no data
# Adds a row for "Adc.c.macros_2.csv" (a CSV that turned out to have no
# data in it), and bumps the macro counts for the CSV scanned after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 ("Adc_Data.c.macros.csv"), which
# shifts it down to row 4 and leaves row 3 free for the new entry.
$ws.Rows(3).Insert()

# Copy the data-row look (wrap-text style) down onto the freshly inserted row.
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)   # xlPasteFormats

# New row for the CSV file that had no data in it.
$ws.Range("A3").Value = "Adc.c.macros_2.csv"

# The counters below are text in this sheet (not real numbers), same as the
# rest of the column -- use TEXT()+paste-values so they land as plain shared
# strings instead of numeric cells / picking up a new number format.
$ws.Range("B3").Formula = '=TEXT(0,"0")'
$ws.Range("C3").Formula = '=TEXT(0,"0")'
$ws.Range("B3:C3").Copy()
$ws.Range("B3:C3").PasteSpecial(-4163)   # xlPasteValues

# Macro counts shift by one for the files scanned after the new (empty) file.
$ws.Range("B2").Formula = '=TEXT(46967,"0")'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("B4").Formula = '=TEXT(2311,"0")'
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)

$excel.CutCopyMode = $false
